# VT_REG1_ELCDMD_K03.xlsx edit script
# Adds an auxiliary "ELC_PRIS" output (electricity price, dummy stand-in for
# export) to the ELC_DMD process: a new commodity row on "Commodities", a new
# OUTPUT~ELC_PRIS column on "Stocks" with data for the ELC_DMD process row,
# and shifts the year columns on "Demand" by one year (2020-2023 -> 2021-2024)
# to make room, following the Commodities row insert that the ELC_DEM
# commodity reference (row 5 -> row 6) now requires.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Commodities sheet: insert a new commodity row for ELC_PRIS above the
#    existing ELC (NRG) row, following the same "NRG" (fuel) pattern used
#    for ELC.
# ---------------------------------------------------------------------
$wsComm = $wb.Worksheets.Item("Commodities")
$wsComm.Rows("5:5").Insert()

$wsComm.Range("B5").Value = "NRG"
$wsComm.Range("D5").Value = "ELC_PRIS"
$wsComm.Range("E5").Value = "Electricity price - Dummy stand in for export"
$wsComm.Range("F5").Value = "PJ"
$wsComm.Range("G5").Value = "FX"
$wsComm.Range("H5").Value = "DAYNITE"

# ---------------------------------------------------------------------
# 2) Stocks sheet: insert a new column before the NCAP_FOM column (L) to
#    hold the OUTPUT~ELC_PRIS auxiliary output, then fill in the header and
#    the ELC_DMD process row's output data.
# ---------------------------------------------------------------------
$wsStocks = $wb.Worksheets.Item("Stocks")
$wsStocks.Columns("L:L").Insert()
$wsStocks.Columns("L:L").ColumnWidth = 15.77734375

$wsStocks.Range("L8").Value = "OUTPUT~ELC_PRIS"

$wsStocks.Range("G9").Value = "ELC_PRIS"
$wsStocks.Range("H9").Value = 2022
$wsStocks.Range("L9").Value = 1

# columns that used to be O/P (the capacity and the DAYNITE factor) shifted
# right by one into P/Q when column L was inserted; O9 is now empty.
$wsStocks.Range("O9").ClearContents()
$wsStocks.Range("P9").Value = 50
$wsStocks.Range("Q9").Formula = "=8760*3.6/10^6"

# A few blank formatted rows below the table (placeholders under the new
# OUTPUT~ELC_PRIS column).
$wsStocks.Range("L13").Value = $wsStocks.Range("L13").Value
$wsStocks.Range("L14").Value = $wsStocks.Range("L14").Value
$wsStocks.Range("L15").Value = $wsStocks.Range("L15").Value

# ---------------------------------------------------------------------
# 3) Demand sheet: the Commodities row insert shifted the "Electricity "
#    commodity description from row 5 to row 6, so the formulas here need
#    to follow; also nudge the year series forward by one year.
# ---------------------------------------------------------------------
$wsDemand = $wb.Worksheets.Item("Demand")
$wsDemand.Range("D5").Formula = "=Commodities!`$E`$6"
$wsDemand.Range("D6").Formula = "=Commodities!`$E`$6"
$wsDemand.Range("D7").Formula = "=Commodities!`$E`$6"
$wsDemand.Range("D8").Formula = "=Commodities!`$E`$6"

$wsDemand.Range("F5").Value = 2021
$wsDemand.Range("F6").Value = 2022
$wsDemand.Range("F7").Value = 2023
$wsDemand.Range("F8").Value = 2024

# ---------------------------------------------------------------------
# 4) Selections / active sheet to mirror the saved UI state.
# ---------------------------------------------------------------------
$wsComm.Range("B5").Select()
$wsProc = $wb.Worksheets.Item("Processes")
$wsProc.Range("D4").Select()
$wsDemand.Range("D16").Select()

$wsStocks.Range("E3").Select()
$wsStocks.Range("H10").Select()

$wsStocks.Activate()
